$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the standalone "Summary" section header row (old row 39). This
#    shifts the final five summary rows (old 40-44) up to become rows 39-43,
#    matching the new dimension A1:B43.
# ---------------------------------------------------------------------------
$ws.Rows.Item(39).Delete()

# ---------------------------------------------------------------------------
# 2. Prefix each branch's "New nominations / Carryover nominations /
#    Confirmed / Unconfirmed / Withdrawn" sub-rows with their section label,
#    e.g. "New nominations" -> "Civilian, New nominations".
# ---------------------------------------------------------------------------
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Carryover nominations"
$ws.Range("A9").Value  = "     Civilian, Confirmed"
$ws.Range("A10").Value = "     Civilian, Unconfirmed "
$ws.Range("A11").Value = "     Civilian, Withdrawn "

$ws.Range("A13").Value = "     Civilian (FS, PHS, CG, NOAA), New nominations"
$ws.Range("A14").Value = "     Civilian (FS, PHS, CG, NOAA), Carryover nominations"
$ws.Range("A15").Value = "     Civilian (FS, PHS, CG, NOAA), Confirmed"
$ws.Range("A16").Value = "     Civilian (FS, PHS, CG, NOAA), Unconfirmed "

$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Carryover nominations"
$ws.Range("A20").Value = "     Air Force, Confirmed"
$ws.Range("A21").Value = "     Air Force, Unconfirmed "
$ws.Range("A22").Value = "     Air Force, Withdrawn "

$ws.Range("A24").Value = "     Army, New nominations"
$ws.Range("A25").Value = "     Army, Carryover nominations"
$ws.Range("A26").Value = "     Army, Confirmed"
$ws.Range("A27").Value = "     Army, Unconfirmed "
$ws.Range("A28").Value = "     Army, Withdrawn "

$ws.Range("A30").Value = "     Navy, New nominations"
$ws.Range("A31").Value = "     Navy, Carryover nominations"
$ws.Range("A32").Value = "     Navy, Confirmed"
$ws.Range("A33").Value = "     Navy, Unconfirmed "
$ws.Range("A34").Value = "     Navy, Withdrawn "

$ws.Range("A36").Value = "     Marine Corps, New nominations"
$ws.Range("A37").Value = "     Marine Corps, Confirmed "
$ws.Range("A38").Value = "     Marine Corps, Unconfirmed "

# ---------------------------------------------------------------------------
# 3. Re-label and swap the final two summary totals (rows 39 & 40 after the
#    row deletion above). The row that used to read "Total carried over from
#    first session" (value 167, General format) becomes "Total carryover
#    nominations", and the row that used to read "Total nominations received
#    this session" (value 31077, #,##0 format) becomes "Total new
#    nominations" - with the two rows trading places. Swap the number
#    formats via a scratch cell so the #,##0 formatting follows the 31077
#    value and the General formatting follows the 167 value, without
#    introducing any new style entries.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1")

$ws.Range("B40").Copy()
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B39").Copy()
$ws.Range("B40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$scratch.Copy()
$ws.Range("B39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$scratch.EntireColumn.Delete()
$excel.CutCopyMode = 0

$ws.Range("A39").Value = "Total new nominations"
$ws.Range("B39").Value = 31077
$ws.Range("A40").Value = "Total carryover nominations"
$ws.Range("B40").Value = 167
